$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91 / column A: tiny re-save precision correction on the existing date-time serial
$ws.Range("A91").Value = 44404.76787926736

# New row 92 - data retrieved Wed Jul 28 18:27:56 UTC 2021
$ws.Range("A92").Value = 44405.76940213077
$ws.Range("B92").Value = 80641
$ws.Range("C92").Value = 68021
$ws.Range("D92").Value = 3754
$ws.Range("E92").Value = 2239
$ws.Range("F92").Value = 1623
$ws.Range("G92").Value = 21182
$ws.Range("H92").Value = 1650
$ws.Range("I92").Value = 913
$ws.Range("J92").Value = 197
